# Update the "想去人数" (wanted-to-go count) figures in column F for the
# "展览" (Exhibition) and "全部类型" (All types) sheets, reflecting the
# regenerated data output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 70
$ws1.Range("F5").Value = 283
$ws1.Range("F6").Value = 393
$ws1.Range("F7").Value = 246
$ws1.Range("F8").Value = 2317
$ws1.Range("F10").Value = 5767
$ws1.Range("F11").Value = 141

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 70
$ws4.Range("F6").Value = 283
$ws4.Range("F7").Value = 393
$ws4.Range("F8").Value = 246
$ws4.Range("F11").Value = 2317
$ws4.Range("F13").Value = 5767
$ws4.Range("F14").Value = 141
